$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add question text to C2 (existing wrapped style s=2 stays as-is)
$ws.Range("C2").Value = ".idea fronete? webstormo bajeriai ar missclickas"

# Row 3: add question text to C3, and grow the row height to fit two lines
$ws.Range("C3").Value = "node_modules pabego 3 moduliai is fronto, webstrom prikolas ar missclickas?"
$ws.Rows.Item(3).RowHeight = 30

# New rows 4-6: additional questions/suggestions, styled like the
# "Explanatory Text" cells but without word-wrap (style index 1)
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)

$ws.Range("C4:C6").WrapText = $false

$ws.Range("C4").Value = "siulau db kataloga kist i doc"
$ws.Range("C5").Value = "paveldejimas, viena didele su daug nulu kur kam nepriklauso laukai ar paveldinti savo lentele? kolkas padaryat paveldinti savo"
$ws.Range("C6").Value = "del medzio, userius siulau viduj palikt tris katalogus nes nebus vieno serviso vieno kontrolerio tiesiog useriui, ar bus"

# Move the active selection the way Excel left it after the edits
$ws.Range("C7").Select()
